$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "reason" -> "tech_reason"
$ws.Range("C1").Value = "tech_reason"

# Row 2 updates
$ws.Range("C2").Value = "The candidate lacks NLP experience, which is a key requirement for the role. Although the candidate has experience in related areas like data analysis and computer vision, the absence of NLP experience affects the overall suitability for the position."
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "The applicant has demonstrated a strong interest in AI and willingness to adapt to a new culture, which aligns with the company's values. The candidate's ability to work well in a team and address challenges through collaboration reflects their suitability for the role. However, the candidate's preference to not work alone may need to be considered within the team dynamics."

# Row 3 updates
$ws.Range("C3").Value = "The candidate has strong experience in NLP, Computer Vision, and PyTorch, which are essential skills for the job. Although the candidate lacks experience in Tensorflow, their projects demonstrate proficiency in related technologies such as image processing, deep learning, and machine learning, making them a strong fit for the position."
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "The candidate's responses indicate a strong interest in AI and a willingness to work in Japan. They demonstrate good communication skills, ability to work in a team, problem-solving skills, and adaptability. The candidate's willingness to learn a new language aligns with the company's values, making them a suitable candidate for the role."
